$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price/volume figures for the symbol list refresh (Mon Feb 6 06:42:55 UTC 2023)
$updates = @(
    @{Cell="D2"; Value="324.86"},
    @{Cell="E2"; Value="-2.43%"},
    @{Cell="D3"; Value="44.62"},
    @{Cell="E3"; Value="1.38%"},
    @{Cell="D4"; Value="5.487"},
    @{Cell="E4"; Value="-5.13%"},
    @{Cell="D5"; Value="0.08040"},
    @{Cell="E5"; Value="-3.68%"},
    @{Cell="E6"; Value="-1.76%"},
    @{Cell="D7"; Value="4.298"},
    @{Cell="E7"; Value="-4.51%"},
    @{Cell="D8"; Value="1.893"},
    @{Cell="E8"; Value="-4.18%"},
    @{Cell="D9"; Value="2.684"},
    @{Cell="E9"; Value="-6.97%"},
    @{Cell="D10"; Value="0.9401"},
    @{Cell="E10"; Value="0.69%"},
    @{Cell="D11"; Value="0.1174"},
    @{Cell="E11"; Value="-5.85%"},
    @{Cell="D12"; Value="0.1865"},
    @{Cell="E12"; Value="-4.25%"},
    @{Cell="D13"; Value="0.09911"},
    @{Cell="E13"; Value="4.24%"},
    @{Cell="D14"; Value="0.04205"},
    @{Cell="E14"; Value="6.77%"},
    @{Cell="E15"; Value="-0.03%"},
    @{Cell="D16"; Value="0.001272"},
    @{Cell="E16"; Value="-2.29%"},
    @{Cell="D17"; Value="0.005864"},
    @{Cell="E17"; Value="-3.57%"},
    @{Cell="D18"; Value="3.592"},
    @{Cell="E18"; Value="2.42%"},
    @{Cell="E19"; Value="-0.71%"},
    @{Cell="D20"; Value="8.553"},
    @{Cell="E20"; Value="-5.66%"},
    @{Cell="E21"; Value="-1.10%"},
    @{Cell="D22"; Value="0.2656"},
    @{Cell="E22"; Value="3.23%"},
    @{Cell="D23"; Value="0.04250"},
    @{Cell="E23"; Value="-3.85%"},
    @{Cell="E24"; Value="-1.56%"},
    @{Cell="D25"; Value="0.004459"},
    @{Cell="E25"; Value="1.71%"},
    @{Cell="E26"; Value="1.06%"},
    @{Cell="D27"; Value="0.0004001"},
    @{Cell="E27"; Value="0.22%"},
    @{Cell="D39"; Value="0.02641"},
    @{Cell="E39"; Value="-6.65%"},
    @{Cell="D40"; Value="0.05492"},
    @{Cell="E40"; Value="-3.76%"},
    @{Cell="D41"; Value="0.007697"},
    @{Cell="E41"; Value="-2.71%"},
    @{Cell="D42"; Value="0.1395"},
    @{Cell="E42"; Value="-2.35%"},
    @{Cell="D43"; Value="0.007328"},
    @{Cell="E43"; Value="-18.99%"},
    @{Cell="E44"; Value="-4.82%"},
    @{Cell="D45"; Value="0.008712"},
    @{Cell="E45"; Value="-14.44%"},
    @{Cell="D46"; Value="0.00007107"},
    @{Cell="E46"; Value="-1.84%"},
    @{Cell="E47"; Value="0.24%"},
    @{Cell="D48"; Value="0.003537"},
    @{Cell="E48"; Value="7.09%"},
    @{Cell="D49"; Value="0.002277"},
    @{Cell="E49"; Value="-0.14%"},
    @{Cell="D50"; Value="0.00002106"},
    @{Cell="E50"; Value="0.24%"},
    @{Cell="D51"; Value="0.0002006"},
    @{Cell="E51"; Value="0.24%"}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (and their
    # significant trailing zeros, e.g. "0.08040") are preserved
    # exactly as authored instead of being coerced to numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
